# ----------------------------------------------------------------------
# shopingFile.xlsx update
# "every funcitons are working perfectly. Now it is time to write the
#  flow chart and report"
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "items": refresh the product catalogue, drop the last 2 rows
# ------------------------------------------------------------------
$items = $wb.Worksheets.Item("items")

$itemsData = @(
    @(1007, "Ipad",        1000, 3),
    @(1008, "mouse",        300, 20),
    @(1009, "Keyboard",     100, 10),
    @(1011, "Mouse",        101, 3),
    @(1013, "Hard disk",   8000, 2),
    @(1014, "SSD",         2800, 25),
    @(1015, "RAM",         4000, 25),
    @(1016, "Monitor",    15000, 3),
    @(1018, "Fiver cable",  100, 998),
    @(1019, "pen drive",   1000, 98),
    @(1020, "Ups",        10000, 4)
)

$r = 5
foreach ($row in $itemsData) {
    $items.Cells.Item($r, 1).Value2 = $row[0]
    $items.Cells.Item($r, 2).Value2 = $row[1]
    $items.Cells.Item($r, 3).Value2 = $row[2]
    $items.Cells.Item($r, 4).Value2 = $row[3]
    $r = $r + 1
}

# rows 16 & 17 no longer exist after the edit
$items.Rows.Item(17).Delete()
$items.Rows.Item(16).Delete()

[void]$items.Range("B9").Select()

# ------------------------------------------------------------------
# Sheet "soldProduct": add a "Time" column header and append new
# sold-product rows (10-21, with an empty separator at row 18)
# ------------------------------------------------------------------
$sold = $wb.Worksheets.Item("soldProduct")

$sold.Range("J1").Value2 = "Time"

$soldData = @(
    @(1015, "RAM",              2, 3500,  7000, "Mohammad",             "fsdaf",    "431",         "bkash"),
    @(1015, "RAM",              3, 3500, 10500, "Mohammad Barkatullah", "fsdf",     "44234",       "COD"),
    @(1019, "Wireless mouse",   3, 1000,  3000, "Mohammad",             "fdsf",     "4234",        "nogod"),
    @(1020, "Bluetooth dangle", 2,  300,   600, "Sadia Afroz",          "uttara",   "34324",       "bkash"),
    @(1018, "Fiver cable",     10,  105,  1050, "Israt Rimpi",          "fdsf",     "4324",        "COD"),
    @(1018, "Fiver cable",     10,  105,  1050, "Partho",               "fsdf",     "41234",       "COD"),
    @(1018, "Fiver cable",     10,  105,  1050, "Partho",               "fsdf",     "54324323425", "nogod"),
    @(1011, "Mouse",            3,  101,   303, "Dipanker",             "fsdf",     "4234",        "bkash")
)

$r = 10
foreach ($row in $soldData) {
    $sold.Cells.Item($r, 1).Value2 = $row[0]
    $sold.Cells.Item($r, 2).Value2 = $row[1]
    $sold.Cells.Item($r, 3).Value2 = $row[2]
    $sold.Cells.Item($r, 4).Value2 = $row[3]
    $sold.Cells.Item($r, 5).Value2 = $row[4]
    $sold.Cells.Item($r, 6).Value2 = $row[5]
    $sold.Cells.Item($r, 7).Value2 = $row[6]
    $sold.Cells.Item($r, 8).Value2 = $row[7]
    $sold.Cells.Item($r, 9).Value2 = $row[8]
    $r = $r + 1
}

# row 18 is a blank separator; only I18 carries the new date/time style
$sold.Range("I18").NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"

$soldData2 = @(
    @(19, 1007, "Ipad", 2,  1000,  2000, "sadia",                 "uttara",   "32443",  "cod"),
    @(20, 1020, "Ups",  2, 10000, 20000, "Mohammad Barkatullah",  "fsdafasd", "435",    "COD"),
    @(21, 1018, "Fiver cable", 2, 100,    200, "Barkat",           "fdsf",     "5443124","COD")
)

foreach ($row in $soldData2) {
    $rr = $row[0]
    $sold.Cells.Item($rr, 1).Value2 = $row[1]
    $sold.Cells.Item($rr, 2).Value2 = $row[2]
    $sold.Cells.Item($rr, 3).Value2 = $row[3]
    $sold.Cells.Item($rr, 4).Value2 = $row[4]
    $sold.Cells.Item($rr, 5).Value2 = $row[5]
    $sold.Cells.Item($rr, 6).Value2 = $row[6]
    $sold.Cells.Item($rr, 7).Value2 = $row[7]
    $sold.Cells.Item($rr, 8).Value2 = $row[8]
    $sold.Cells.Item($rr, 9).Value2 = $row[9]
}

[void]$sold.Range("F10").Select()

# ------------------------------------------------------------------
# Sheet "userAccount": rewrite the account table contents
# ------------------------------------------------------------------
$users = $wb.Worksheets.Item("userAccount")

$users.Range("A2").Value2 = "Mohammad "
$users.Range("B2").Value2 = "barkat"
$users.Range("C2").Value2 = "123"
$users.Range("D2").Value2 = "fsdf"
$users.Range("E2").Value2 = "4324"

$users.Range("A3").Value2 = "Barkat"
$users.Range("B3").Value2 = "barkatopu"
$users.Range("C3").Value2 = "123"
$users.Range("D3").Value2 = "fdsf"
$users.Range("E3").Value2 = "4324"

$users.Range("A4").Value2 = "Dip"
$users.Range("B4").Value2 = "dip"
$users.Range("C4").Value2 = "123"
$users.Range("D4").Value2 = "vasfd"
$users.Range("E4").Value2 = "4234"

$users.Range("A5").Value2 = "partho"
$users.Range("B5").Value2 = "partho"
$users.Range("C5").Value2 = "123"
$users.Range("D5").Value2 = "fsdf"
$users.Range("E5").Value2 = "4324"

$users.Range("A6").Value2 = "Mohammad Barkatullah"
$users.Range("B6").Value2 = "barkat1345"
$users.Range("C6").Value2 = "1234"
$users.Range("D6").Value2 = "xyz"
$users.Range("E6").Value2 = "01521206720"

$users.Range("A7").Value2 = "Sadia afroz"
$users.Range("B7").Value2 = "sadia"
$users.Range("C7").Value2 = "123"
$users.Range("D7").Value2 = "asd"
$users.Range("E7").Value2 = "543543"

$users.Range("A8").Value2 = "Israt rimpi"
$users.Range("B8").Value2 = "rimpi"
$users.Range("C8").Value2 = "123"
$users.Range("D8").Value2 = "uttara"
$users.Range("E8").Value2 = "341234"

[void]$users.Range("E16").Select()

# ------------------------------------------------------------------
# Make "soldProduct" the active sheet/tab (workbook activeTab goes
# from items-context 2 to 1)
# ------------------------------------------------------------------
$sold.Activate()
[void]$sold.Range("F10").Select()
